$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Código Barras" column (column C) entirely - the remaining
# headers (Laboratorio, N° Lote, Fecha Vencimiento, Precio Venta, Cantidad),
# their cell formatting, and the custom widths on the old F/G columns all
# shift one column to the left along with it.
$ws.Range("C1").EntireColumn.Delete() | Out-Null

# Update the active selection to match the post-edit state.
$ws.Range("E5").Select() | Out-Null
